$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 11.6861
$ws.Range("E4").Value = 13.73869999999999
$ws.Range("E7").Value = 12.1143
$ws.Range("E8").Value = 12.1949
$ws.Range("D11").Value = -8.5108
$ws.Range("D12").Value = -8.477499999999999
$ws.Range("E12").Value = 12.2815
$ws.Range("E14").Value = 14.0776
$ws.Range("D15").Value = -7.975700000000001
$ws.Range("E22").Value = 11.76
